$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Append a new row of data (row 43) to the kaspa_buys sheet, reflecting
# the 2025-07-20 run: Date, Coins, Price, Cost.

# Use a leading apostrophe so Excel stores the date-like text as a literal
# string instead of auto-converting it to a date serial number, then reset
# the cell style to Normal so it matches the other text-date rows (no
# custom number format is left behind).
$ws.Range("A43").Value = "'07/20/2025"
$ws.Range("A43").Style = "Normal"

$ws.Range("B43").Value = 100.8229999999967
$ws.Range("C43").Value = 0.09918371800085625
$ws.Range("D43").Value = 10
